# Append new product rows to the price list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(92007895, "Jogo de Brocas Talhadeira e Ponteiro para Concreto SDS 6 Peças Dexter", 49.9),
    @(92007762, "Jogo de Brocas com 5 para Madeira 6 para Metal e 5 para Alvenaria 16 Peças Dexter", 59.9),
    @(89833870, "Jogo de Brocas de Vídea encaixe SDS Plus para Concreto 6 a 12mm 5 Peças Dexter", 58.9),
    @(88162186, "Lixeira para Banheiro 3L Metal Aço Pedal Prata Escovado Sensea", 52.9),
    @(89408333, "Sacos de Lixo Azul 50L Contém 50 Sacos Embalixo", 19.79),
    @(89408270, "Sacos de Lixo Branco Pia e Banheiro Contém 100 Sacos 38x40cm Embalixo", 19.79),
    @(90672694, "Saco Lixo Embalixo 100L Wecycle Reforçado Azul Com 10 Sacos", 19.79),
    @(90672680, "Saco Lixo Embalixo 50L Reforçado Preto Com 50 Sacos", 59.9),
    @(88187246, "Cabo de Aluminio 1,50m Bettanin", 52.9),
    @(90351765, "Martelete Perfurador Rompedor Makita com Maleta com 5 Brocas e Acessórios SDS Plus 800W 2,7J HR 2470 220V", 769.9)
)

# Find the last used row so the new rows are appended right after
# the existing data (row 127 -> new rows start at 128).
$lastRow = $ws.UsedRange.Rows.Count
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
